$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The first three "Newuserstaging_1/2/3" login rows were bad/duplicate
# journey-page references. Remove rows 2-4 entirely so every following
# row (previously row 5 onward) shifts up by three.
$ws.Rows("2:4").Delete()
